# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '66.027.56'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = '2.691.25'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextCell 'D5' '608.59'
$ws.Range('E5').Value = '  +1.10%  '
Set-TextCell 'D6' '158.09'
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +2.87%  '
Set-TextCell 'D10' '5.99'
$ws.Range('E10').Value = '  +3.33%  '
$ws.Range('E11').Value = '  -2.95%  '
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('E13').Value = '  +8.31%  '
Set-TextCell 'D14' '30.02'
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('D15').Value = '3.172.50'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').Value = '65.811.85'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = '2.688.16'
$ws.Range('E17').Value = '  +2.12%  '
Set-TextCell 'D18' '12.74'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('E19').Value = '  -1.23%  '
Set-TextCell 'D20' '7.71'
$ws.Range('E20').Value = '  +4.27%  '
Set-TextCell 'D21' '356.90'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('E22').Value = '  +0.24%  '
Set-TextCell 'D23' '71.01'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('E24').Value = '  +15.99%  '
Set-TextCell 'D25' '9.97'
$ws.Range('E25').Value = '  +5.63%  '
Set-TextCell 'D26' '1.62'
$ws.Range('E26').Value = '  -5.37%  '
$ws.Range('E27').Value = '  -0.09%  '
Set-TextCell 'D28' '0.172'
$ws.Range('E28').Value = '  +3.62%  '
Set-TextCell 'D29' '8.25'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('E31').Value = '  -0.24%  '
Set-TextCell 'D32' '531.95'
$ws.Range('E32').Value = '  -3.81%  '
$ws.Range('E33').Value = '  -2.08%  '
Set-TextCell 'D34' '6.62'
$ws.Range('E34').Value = '  +3.73%  '
Set-TextCell 'D35' '5.48'
$ws.Range('E35').Value = '  -2.98%  '
Set-TextCell 'D36' '0.432'
$ws.Range('E36').Value = '  -0.07%  '
Set-TextCell 'D37' '20.70'
$ws.Range('E37').Value = '  +0.92%  '
Set-TextCell 'D38' '161.75'
$ws.Range('E38').Value = '  -0.25%  '
Set-TextCell 'D39' '2.00'
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D42' '42.36'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D43' '167.53'
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('E44').Value = '  -1.51%  '
Set-TextCell 'D45' '0.0632'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('E50').Value = '  +5.38%  '
$ws.Range('E51').Value = '  +1.37%  '
